$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 37, shifting current rows 37-74 down to 38-75
$ws.Rows("37").Insert()

# Populate the newly inserted row 37 with data (same as the row that is now 38,
# i.e. the original row 37, except for the Fecha (D) and Volumen (J) values).
$ws.Range("A37").Value = 7
$ws.Range("B37").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C37").Value = "Ñuble"
$ws.Range("D37").NumberFormat = $ws.Range("D38").NumberFormat
$ws.Range("D37").Value = 44580
$ws.Range("E37").Value = 16
$ws.Range("F37").Value = 100112030
$ws.Range("G37").Value = "Poroto granado"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 23000
$ws.Range("L37").Value = 24000
$ws.Range("M37").Value = 23500
$ws.Range("N37").Value = "`$/saco 25 kilos"
$ws.Range("O37").Value = "Provincia de Diguillín"
$ws.Range("P37").Value = 940
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
